$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("BD2").Value = 126
$ws.Range("G3").Value = 1.73
$ws.Range("I3").Value = 5.75
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 6
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("AF3").Value = 81
$ws.Range("AJ3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AV3").Value = 81
$ws.Range("AX3").Value = 34
$ws.Range("BB3").Value = 451
$ws.Range("I9").Value = 2.63
$ws.Range("L9").Value = 3.4
$ws.Range("Z9").Value = 29
$ws.Range("AI9").Value = 10
$ws.Range("AN9").Value = 4.75
$ws.Range("AO9").Value = 17
